$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 744.9091
$ws.Range("I18").Value = 759.4
$ws.Range("J18").Value = 600
$ws.Range("K18").Value = 759.4
$ws.Range("L18").Value = 600
$ws.Range("M18").Value = -475.4
$ws.Range("N18").Value = -1168
$ws.Range("H32").Value = 6920.154
$ws.Range("I32").Value = 7483.5
$ws.Range("J32").Value = 6817.727
$ws.Range("K32").Value = 7483.5
$ws.Range("L32").Value = 6817.727
$ws.Range("M32").Value = -7157.5
$ws.Range("N32").Value = -7469.727
$ws.Range("H52").Value = 5261.2
$ws.Range("I52").Value = 3768.6667
$ws.Range("K52").Value = 11306.0001
$ws.Range("M52").Value = -11146.0001
$ws.Range("H82").Value = 8334376.5
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 8334376.5
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H116").Value = 6235.077
$ws.Range("I116").Value = 6999
$ws.Range("J116").Value = 6171.4165
$ws.Range("K116").Value = 6999
$ws.Range("L116").Value = 6171.4165
$ws.Range("M116").Value = -3557
$ws.Range("N116").Value = -13055.4165
$ws.Range("H138").Value = 1971.3
$ws.Range("J138").Value = 2141.3333
$ws.Range("L138").Value = 6423.999899999999
$ws.Range("N138").Value = -16703.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3100.3845
$ws.Range("I45").Value = 3499.25
$ws.Range("J45").Value = 2923.111
$ws.Range("K45").Value = 3499.25
$ws.Range("L45").Value = 2923.111
$ws.Range("M45").Value = -3122.25
$ws.Range("N45").Value = -3677.111
$ws.Range("H74").Value = 2580.8823
$ws.Range("I74").Value = 2132.9092
$ws.Range("K74").Value = 2132.9092
$ws.Range("M74").Value = -1258.9092
$ws.Range("H77").Value = 2580.8823
$ws.Range("I77").Value = 2132.9092
$ws.Range("K77").Value = 10664.546
$ws.Range("M77").Value = -6296.546
$ws.Range("H97").Value = 450.5
$ws.Range("I97").Value = 373.33334
$ws.Range("J97").Value = 604.8333
$ws.Range("K97").Value = 373.33334
$ws.Range("L97").Value = 604.8333
$ws.Range("M97").Value = 122.66666
$ws.Range("N97").Value = -1596.8333
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1805.3125
$ws.Range("I94").Value = 907.25
$ws.Range("J94").Value = 4499.5
$ws.Range("K94").Value = 907.25
$ws.Range("L94").Value = 4499.5
$ws.Range("M94").Value = -456.25
$ws.Range("N94").Value = -5401.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5543
$ws.Range("I16").Value = 6153.875
$ws.Range("J16").Value = 4321.25
$ws.Range("K16").Value = 6153.875
$ws.Range("L16").Value = 4321.25
$ws.Range("M16").Value = -5866.875
$ws.Range("N16").Value = -4895.25
$ws.Range("H25").Value = 20000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 20000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 20000
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -20348
$ws.Range("H98").Value = 26999.666
$ws.Range("J98").Value = 26999.666
$ws.Range("L98").Value = 26999.666
$ws.Range("N98").Value = -31491.666
$ws.Range("H99").Value = 88891330
$ws.Range("I99").Value = 66668932
$ws.Range("J99").Value = 111113730
$ws.Range("K99").Value = 66668932
$ws.Range("L99").Value = 111113730
$ws.Range("M99").Value = -66667434
$ws.Range("N99").Value = -111116726
$ws.Range("H113").Value = 5543
$ws.Range("I113").Value = 6153.875
$ws.Range("J113").Value = 4321.25
$ws.Range("K113").Value = 6153.875
$ws.Range("L113").Value = 4321.25
$ws.Range("M113").Value = -3983.875
$ws.Range("N113").Value = -8661.25
$ws.Range("H126").Value = 88891330
$ws.Range("I126").Value = 66668932
$ws.Range("J126").Value = 111113730
$ws.Range("K126").Value = 200006796
$ws.Range("L126").Value = 333341190
$ws.Range("M126").Value = -200004326
$ws.Range("N126").Value = -333346130

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 392.44446
$ws.Range("I12").Value = 62.11111
$ws.Range("J12").Value = 722.7778
$ws.Range("K12").Value = 186.33333
$ws.Range("L12").Value = 2168.3334
$ws.Range("M12").Value = -13.33332999999999
$ws.Range("N12").Value = -2514.3334
$ws.Range("H22").Value = 867.7273
$ws.Range("I22").Value = 725
$ws.Range("J22").Value = 882
$ws.Range("K22").Value = 2175
$ws.Range("L22").Value = 2646
$ws.Range("M22").Value = -2006
$ws.Range("N22").Value = -2984
$ws.Range("H27").Value = 867.7273
$ws.Range("I27").Value = 725
$ws.Range("J27").Value = 882
$ws.Range("K27").Value = 2175
$ws.Range("L27").Value = 2646
$ws.Range("M27").Value = -2073
$ws.Range("N27").Value = -2850
$ws.Range("H41").Value = 3460.2
$ws.Range("J41").Value = 3460.2
$ws.Range("L41").Value = 10380.6
$ws.Range("N41").Value = -11056.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 9579.1
$ws.Range("I22").Value = 25498.666
$ws.Range("J22").Value = 2756.4285
$ws.Range("K22").Value = 25498.666
$ws.Range("L22").Value = 2756.4285
$ws.Range("M22").Value = -24969.666
$ws.Range("N22").Value = -3814.4285
$ws.Range("H102").Value = 4528.0454
$ws.Range("I102").Value = 4528.0454
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4528.0454
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2906.0454
$ws.Range("N102").ClearContents()
$ws.Range("H126").Value = 5586.5557
$ws.Range("J126").Value = 5891.9
$ws.Range("L126").Value = 17675.7
$ws.Range("N126").Value = -22615.7
$ws.Range("H132").Value = 2872.5
$ws.Range("I132").Value = 2864.6667
$ws.Range("K132").Value = 8594.000100000001
$ws.Range("M132").Value = -6064.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1736.4
$ws.Range("J46").Value = 1697.5
$ws.Range("L46").Value = 1697.5
$ws.Range("N46").Value = -2073.5
$ws.Range("H55").Value = 362.18182
$ws.Range("I55").Value = 388.25
$ws.Range("J55").Value = 292.66666
$ws.Range("K55").Value = 388.25
$ws.Range("L55").Value = 292.66666
$ws.Range("M55").Value = -215.25
$ws.Range("N55").Value = -638.66666
$ws.Range("H132").Value = 2583.1794
$ws.Range("I132").Value = 2353.1292
$ws.Range("K132").Value = 7059.3876
$ws.Range("M132").Value = -4529.3876
$ws.Range("H136").Value = 3399
$ws.Range("I136").Value = 3523.7
$ws.Range("K136").Value = 10571.1
$ws.Range("M136").Value = -8021.099999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1698080.2
$ws.Range("J96").Value = 1014.5455
$ws.Range("L96").Value = 1014.5455
$ws.Range("N96").Value = -3760.5455
$ws.Range("H113").Value = 1054
$ws.Range("I113").Value = 975.44446
$ws.Range("K113").Value = 2926.33338
$ws.Range("M113").Value = -756.33338
$ws.Range("H126").Value = 43480210
$ws.Range("I126").Value = 2036.9546
$ws.Range("K126").Value = 6110.8638
$ws.Range("M126").Value = -3640.8638
$ws.Range("H136").Value = 3484.6843
$ws.Range("I136").Value = 2549.6
$ws.Range("K136").Value = 7648.799999999999
$ws.Range("M136").Value = -5098.799999999999
